$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '66.585.68'
$ws.Cells.Item(2, 5).Value = '  +1.74%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.693.86'
$ws.Cells.Item(3, 5).Value = '  +4.61%  '

$ws.Cells.Item(4, 5).Value = '  +0.46%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '419.50'
$ws.Cells.Item(5, 5).Value = '  -0.83%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '130.01'
$ws.Cells.Item(6, 5).Value = '  -2.12%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.686.62'
$ws.Cells.Item(7, 5).Value = '  +4.68%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.644'
$ws.Cells.Item(8, 5).Value = '  +0.46%  '

$ws.Cells.Item(9, 5).Value = '  +0.05%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.768'
$ws.Cells.Item(10, 5).Value = '  -3.28%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.182'
$ws.Cells.Item(11, 5).Value = '  +8.71%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0000397'
$ws.Cells.Item(12, 5).Value = '  +46.64%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '43.14'
$ws.Cells.Item(13, 5).Value = '  -0.57%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '10.64'
$ws.Cells.Item(14, 5).Value = '  +5.55%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.286.04'
$ws.Cells.Item(15, 5).Value = '  +5.13%  '

$ws.Cells.Item(16, 5).Value = '  -0.88%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '20.60'
$ws.Cells.Item(17, 5).Value = '  -0.71%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '3.683.29'
$ws.Cells.Item(18, 5).Value = '  +4.76%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '13.27'
$ws.Cells.Item(19, 5).Value = '  +5.52%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '1.12'
$ws.Cells.Item(20, 5).Value = '  +1.33%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '66.780.94'
$ws.Cells.Item(21, 5).Value = '  +2.63%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '445.15'
$ws.Cells.Item(22, 5).Value = '  -2.33%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '16.50'
$ws.Cells.Item(23, 5).Value = '  +22.95%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '89.95'
$ws.Cells.Item(24, 5).Value = '  -2.14%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.15'
$ws.Cells.Item(25, 5).Value = '  -2.97%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '37.43'
$ws.Cells.Item(26, 5).Value = '  +8.60%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.27'
$ws.Cells.Item(27, 5).Value = '  +0.89%  '

$ws.Cells.Item(28, 5).Value = '  -1.34%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '5.08'
$ws.Cells.Item(29, 5).Value = '  +5.49%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.125'
$ws.Cells.Item(30, 5).Value = '  +9.26%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '12.74'
$ws.Cells.Item(31, 5).Value = '  +0.82%  '

$ws.Cells.Item(32, 5).Value = '  -2.31%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '7.32'
$ws.Cells.Item(33, 5).Value = '  -3.18%  '

$ws.Cells.Item(34, 5).Value = '  +1.96%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '41.45'
$ws.Cells.Item(35, 5).Value = '  +2.82%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '57.37'
$ws.Cells.Item(36, 5).Value = '  -0.52%  '

$ws.Cells.Item(37, 5).Value = '  -0.05%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0494'
$ws.Cells.Item(38, 5).Value = '  -2.87%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0₃0738'
$ws.Cells.Item(39, 5).Value = '  -1.14%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.07'
$ws.Cells.Item(40, 5).Value = '  +32.58%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.151'
$ws.Cells.Item(41, 5).Value = '  +3.89%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '29.49'
$ws.Cells.Item(42, 5).Value = '  +33.28%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.998'
$ws.Cells.Item(43, 5).Value = '  +0.01%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '3.42'
$ws.Cells.Item(44, 5).Value = '  +2.15%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '148.84'
$ws.Cells.Item(45, 5).Value = '  +1.93%  '

$ws.Cells.Item(46, 5).Value = '  +4.19%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).Value = '2.68'
$ws.Cells.Item(47, 5).Value = '  -4.36%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).Value = '2.90'
$ws.Cells.Item(48, 5).Value = '  -7.22%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '4.34'
$ws.Cells.Item(49, 5).Value = '  -5.17%  '

$ws.Cells.Item(50, 5).Value = '  -3.41%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.160'
$ws.Cells.Item(51, 5).Value = '  +11.50%  '
